# =====================================================================
# osticket/scp/investigacion.docx — "PDF Y WORD listo"
#
# 1. Remove the three empty-placeholder paragraphs "Título del
#    proyecto:", "Línea temática afín:" and "Origen del proyecto:"
#    (the form no longer collects those fields).
# 2. The last (previously blank) run of the "Nombre del Autor(es):"
#    paragraph now carries the {usuario} merge field.
# 3. The <w:lastRenderedPageBreak/> bookmark that Word leaves behind
#    from the last render moves forward one placeholder in two spots:
#    {p6} -> {p7} and {p19} -> {p21}.
# =====================================================================

$d = $word.ActiveDocument

$xmlWrapperOpen = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$xmlWrapperClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------
# 1) Delete the "Título del proyecto:" / "Línea temática afín:" /
#    "Origen del proyecto:" paragraphs in full (text + mark).
#    Locate all three first, then delete back-to-front so earlier
#    Range objects stay valid while later ones are removed.
# ---------------------------------------------------------------
$rTitulo = $d.Content
$rTitulo.Find.Execute("Título del proyecto:")
$pTitulo = $rTitulo.Paragraphs(1)

$rLinea = $d.Content
$rLinea.Find.Execute("Línea temática afín")
$pLinea = $rLinea.Paragraphs(1)

$rOrigen = $d.Content
$rOrigen.Find.Execute("Origen del proyecto:")
$pOrigen = $rOrigen.Paragraphs(1)

$pOrigen.Range.Delete()
$pLinea.Range.Delete()
$pTitulo.Range.Delete()

# ---------------------------------------------------------------
# 2) "Nombre del Autor(es): " — the paragraph's final run (a lone
#    space) becomes the {usuario} field; the run ahead of it (also a
#    lone space) is left exactly as it was.
# ---------------------------------------------------------------
$rAutor = $d.Content
$rAutor.Find.Execute("Nombre del Autor(es):")
$pAutor = $rAutor.Paragraphs(1)
$autorEnd = $pAutor.Range.End
# last character before the paragraph mark = the final run's text
$lastRun = $d.Range($autorEnd - 2, $autorEnd - 1)
$lastRun.InsertXML($xmlWrapperOpen + '<w:body><w:p><w:r w:rsidRPr="002A6D4E"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:bCs/><w:sz w:val="20"/></w:rPr><w:t>{usuario}</w:t></w:r></w:p></w:body>' + $xmlWrapperClose)

# ---------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from the {p6} run to the {p7}
#    run, and from the {p19} run to the {p21} run. Each target
#    paragraph is rewritten in place (same paraId/rsid/pPr/rPr) with
#    the break element added or dropped, which is what Word itself
#    does when repagination shifts where the break last fell.
# ---------------------------------------------------------------

# --- {p6}: drop the page-break mark ---
$rP6 = $d.Content
$rP6.Find.Execute("{p6}")
$rP6.Paragraphs(1).Range.InsertXML($xmlWrapperOpen + '<w:body><w:p w14:paraId="3966D285" w14:textId="1B31F8C3" w:rsidR="009567C1" w:rsidRPr="00AA234E" w:rsidRDefault="00970EFB" w:rsidP="001F0712"><w:pPr><w:pStyle w:val="Textoindependiente2"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>{p6}</w:t></w:r></w:p></w:body>' + $xmlWrapperClose)

# --- {p7}: add the page-break mark ---
$rP7 = $d.Content
$rP7.Find.Execute("{p7}")
$rP7.Paragraphs(1).Range.InsertXML($xmlWrapperOpen + '<w:body><w:p w14:paraId="3966D28D" w14:textId="4DC9B2A2" w:rsidR="009567C1" w:rsidRPr="00AA234E" w:rsidRDefault="00970EFB" w:rsidP="001F0712"><w:pPr><w:pStyle w:val="Textoindependiente2"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>{p7}</w:t></w:r></w:p></w:body>' + $xmlWrapperClose)

# --- {p19}: drop the page-break mark ---
$rP19 = $d.Content
$rP19.Find.Execute("{p19}")
$rP19.Paragraphs(1).Range.InsertXML($xmlWrapperOpen + '<w:body><w:p w14:paraId="3966D2F6" w14:textId="435464F0" w:rsidR="009567C1" w:rsidRPr="00A15A72" w:rsidRDefault="00970EFB" w:rsidP="00D912A0"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr><w:t>{p19}</w:t></w:r></w:p></w:body>' + $xmlWrapperClose)

# --- {p21}: add the page-break mark ---
$rP21 = $d.Content
$rP21.Find.Execute("{p21}")
$rP21.Paragraphs(1).Range.InsertXML($xmlWrapperOpen + '<w:body><w:p w14:paraId="3966D306" w14:textId="5CBC02E7" w:rsidR="009567C1" w:rsidRPr="00A15A72" w:rsidRDefault="00970EFB" w:rsidP="00D912A0"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr><w:lastRenderedPageBreak/><w:t>{p21}</w:t></w:r></w:p></w:body>' + $xmlWrapperClose)

Write-Host "Edit complete."
